$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2488856.8
$ws.Range("I33").Value = 3353435.5
$ws.Range("K33").Value = 3353435.5
$ws.Range("M33").Value = -3353206.5

$ws.Range("H98").Value = 5269.0435
$ws.Range("I98").Value = 4431.923
$ws.Range("J98").Value = 6357.3
$ws.Range("K98").Value = 4431.923
$ws.Range("L98").Value = 6357.3
$ws.Range("M98").Value = -2933.923
$ws.Range("N98").Value = -9353.299999999999

$ws.Range("H112").Value = 3986
$ws.Range("J112").Value = 4253.1304
$ws.Range("L112").Value = 12759.3912
$ws.Range("N112").Value = -14975.3912

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H122").Value = 5269.0435
$ws.Range("I122").Value = 4431.923
$ws.Range("J122").Value = 6357.3
$ws.Range("K122").Value = 13295.769
$ws.Range("L122").Value = 19071.9
$ws.Range("M122").Value = -10845.769
$ws.Range("N122").Value = -23971.9

$ws.Range("H125").Value = 200003840
$ws.Range("I125").Value = 500002100
$ws.Range("J125").Value = 4999
$ws.Range("K125").Value = 4500018900
$ws.Range("L125").Value = 44991
$ws.Range("M125").Value = -4500016440
$ws.Range("N125").Value = -49911

$ws.Range("H135").Value = 857.7646999999999
$ws.Range("I135").Value = 942.86664
$ws.Range("J135").Value = 219.5
$ws.Range("K135").Value = 8485.79976
$ws.Range("L135").Value = 1975.5
$ws.Range("M135").Value = -5950.79976
$ws.Range("N135").Value = -7045.5

$ws.Range("H138").Value = 4623.159
$ws.Range("I138").Value = 7058.7393
$ws.Range("J138").Value = 3405.3696
$ws.Range("K138").Value = 21176.2179
$ws.Range("L138").Value = 10216.1088
$ws.Range("M138").Value = -16036.2179
$ws.Range("N138").Value = -20496.1088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 144014.34
$ws.Range("I32").Value = 152090.56
$ws.Range("J32").Value = 10756.75
$ws.Range("K32").Value = 152090.56
$ws.Range("L32").Value = 10756.75
$ws.Range("M32").Value = -151803.56
$ws.Range("N32").Value = -11330.75

$ws.Range("H61").Value = 5885488
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5885488
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 5885488
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -5885912

$ws.Range("H63").Value = 17466.08
$ws.Range("I63").Value = 3412.5557
$ws.Range("J63").Value = 25371.188
$ws.Range("K63").Value = 3412.5557
$ws.Range("L63").Value = 25371.188
$ws.Range("M63").Value = -2726.5557
$ws.Range("N63").Value = -26743.188

$ws.Range("H66").Value = 17466.08
$ws.Range("I66").Value = 3412.5557
$ws.Range("J66").Value = 25371.188
$ws.Range("K66").Value = 17062.7785
$ws.Range("L66").Value = 126855.94
$ws.Range("M66").Value = -13630.7785
$ws.Range("N66").Value = -133719.94

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 1173237.6
$ws.Range("I74").Value = 2474.6667
$ws.Range("J74").Value = 3180259.8
$ws.Range("K74").Value = 2474.6667
$ws.Range("L74").Value = 3180259.8
$ws.Range("M74").Value = -1600.6667
$ws.Range("N74").Value = -3182007.8

$ws.Range("H77").Value = 1173237.6
$ws.Range("I77").Value = 2474.6667
$ws.Range("J77").Value = 3180259.8
$ws.Range("K77").Value = 12373.3335
$ws.Range("L77").Value = 15901299
$ws.Range("M77").Value = -8005.333500000001
$ws.Range("N77").Value = -15910035

$ws.Range("H122").Value = 3723.75
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3723.75
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11171.25
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16071.25

$ws.Range("H132").Value = 3479.6667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3479.6667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 10439.0001
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -15499.0001

$ws.Range("H136").Value = 5885488
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5885488
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 17656464
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17661564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 16998
$ws.Range("I96").Value = 16998
$ws.Range("K96").Value = 16998
$ws.Range("M96").Value = -14252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 2373.5
$ws.Range("I39").Value = 2498.3333
$ws.Range("K39").Value = 2498.3333
$ws.Range("M39").Value = -2107.3333

$ws.Range("H49").Value = 2373.5
$ws.Range("I49").Value = 2498.3333
$ws.Range("K49").Value = 2498.3333
$ws.Range("M49").Value = -2316.3333

$ws.Range("H86").Value = 11588.421
$ws.Range("I86").Value = 7163.091
$ws.Range("J86").Value = 17673.25
$ws.Range("K86").Value = 7163.091
$ws.Range("L86").Value = 17673.25
$ws.Range("M86").Value = -6040.091
$ws.Range("N86").Value = -19919.25

$ws.Range("H89").Value = 11588.421
$ws.Range("I89").Value = 7163.091
$ws.Range("J89").Value = 17673.25
$ws.Range("K89").Value = 35815.455
$ws.Range("L89").Value = 88366.25
$ws.Range("M89").Value = -30199.455
$ws.Range("N89").Value = -99598.25

$ws.Range("H97").Value = 65619.8
$ws.Range("J97").Value = 65619.8
$ws.Range("L97").Value = 65619.8
$ws.Range("N97").Value = -67601.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1793
$ws.Range("I8").Value = 1793
$ws.Range("K8").Value = 5379
$ws.Range("M8").Value = -5240

$ws.Range("H80").Value = 1200
$ws.Range("I80").Value = 1200
$ws.Range("K80").Value = 3600
$ws.Range("M80").Value = -2664

$ws.Range("H83").Value = 1200
$ws.Range("I83").Value = 1200
$ws.Range("K83").Value = 10800
$ws.Range("M83").Value = -6120

$ws.Range("H86").Value = 1937.25
$ws.Range("J86").Value = 1571.1428
$ws.Range("L86").Value = 4713.428400000001
$ws.Range("N86").Value = -7085.428400000001

$ws.Range("H89").Value = 1937.25
$ws.Range("J89").Value = 1571.1428
$ws.Range("L89").Value = 14140.2852
$ws.Range("N89").Value = -25996.2852

$ws.Range("H104").Value = 5924.6
$ws.Range("I104").Value = 2591.3333
$ws.Range("J104").Value = 6757.9165
$ws.Range("K104").Value = 7773.999899999999
$ws.Range("L104").Value = 20273.7495
$ws.Range("M104").Value = -5152.999899999999
$ws.Range("N104").Value = -25515.7495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31251166
$ws.Range("I102").Value = 41667930
$ws.Range("J102").Value = 875
$ws.Range("K102").Value = 41667930
$ws.Range("L102").Value = 875
$ws.Range("M102").Value = -41666308
$ws.Range("N102").Value = -4119

$ws.Range("H122").Value = 2493.8215
$ws.Range("I122").Value = 2249.2222
$ws.Range("J122").Value = 2934.1
$ws.Range("K122").Value = 6747.6666
$ws.Range("L122").Value = 8802.299999999999
$ws.Range("M122").Value = -4297.6666
$ws.Range("N122").Value = -13702.3

$ws.Range("H132").Value = 646063.0600000001
$ws.Range("I132").Value = 5248.871
$ws.Range("J132").Value = 2853312
$ws.Range("K132").Value = 15746.613
$ws.Range("L132").Value = 8559936
$ws.Range("M132").Value = -13216.613
$ws.Range("N132").Value = -8564996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2718.342
$ws.Range("I40").Value = 1537.4
$ws.Range("K40").Value = 1537.4
$ws.Range("M40").Value = -1401.4

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H93").Value = 3024.4614
$ws.Range("I93").Value = 2302.1428
$ws.Range("J93").Value = 3867.1667
$ws.Range("K93").Value = 2302.1428
$ws.Range("L93").Value = 3867.1667
$ws.Range("M93").Value = -1054.1428
$ws.Range("N93").Value = -6363.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1617.76
$ws.Range("I122").Value = 1438.9546
$ws.Range("J122").Value = 2929
$ws.Range("K122").Value = 4316.8638
$ws.Range("L122").Value = 8787
$ws.Range("M122").Value = -1866.8638
$ws.Range("N122").Value = -13687

$ws.Range("H123").Value = 102993.336
$ws.Range("J123").Value = 102993.336
$ws.Range("L123").Value = 102993.336
$ws.Range("N123").Value = -112793.336
